$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 283
$ws.Range("I2").Value = 760
$ws.Range("J2").Value = 3244
$ws.Range("L2").Value = 843
$ws.Range("M2").Value = 57
$ws.Range("N2").Value = 541
$ws.Range("P2").Value = 12
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 35
$ws.Range("S2").Value = 332
$ws.Range("T2").Value = 573
$ws.Range("U2").Value = 44
$ws.Range("V2").Value = 4837
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 4935
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 73
$ws.Range("AA2").Value = 24
